# Add a new time-log entry (row 14): 10/19/2023, Internship,
# "Completed daily operations, 8 hours" entry (reusing existing shared strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column (A) - reuse the date number format / font / alignment from the row above
$ws.Range("A14").Value = 45218
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("A14").Font.Size = $ws.Range("A13").Font.Size
$ws.Range("A14").HorizontalAlignment = $ws.Range("A13").HorizontalAlignment

# Name of Task column (B) - "Internship"
$ws.Range("B14").Value = "Internship"
$ws.Range("B14").Font.Size = $ws.Range("B13").Font.Size

# Description column (C) - reuse existing description text
$ws.Range("C14").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"
$ws.Range("C14").Font.Size = $ws.Range("C13").Font.Size

# Move the active selection to the next empty row, as Excel would after data entry
$ws.Range("C15").Select() | Out-Null
